$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '70.306.14'
$ws.Range('E2').Value = '  -0.98%  '

# Row 3
$ws.Range('D3').Value = '3.768.81'
$ws.Range('E3').Value = '  -1.55%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '693.13'
$ws.Range('E5').Value = '  -1.65%  '

# Row 6
$ws.Range('E6').Value = '  -2.52%  '

# Row 7
$ws.Range('D7').Value = '3.769.57'
$ws.Range('E7').Value = '  -1.50%  '

# Row 8
$ws.Range('E8').Value = '  +0.12%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  -1.19%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.159'
$ws.Range('E10').Value = '  -1.97%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.48'
$ws.Range('E11').Value = '  +0.84%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.473'
$ws.Range('E12').Value = '  +3.10%  '

# Row 13
$ws.Range('E13').Value = '  -3.48%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.64'
$ws.Range('E14').Value = '  -2.56%  '

# Row 16
$ws.Range('E16').Value = '  +0.62%  '

# Row 17
$ws.Range('D17').Value = '70.376.40'
$ws.Range('E17').Value = '  -0.93%  '

# Row 18
$ws.Range('E18').Value = '  -0.02%  '

# Row 19
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.09'
$ws.Range('E19').Value = '  -1.93%  '

# Row 20
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.32'
$ws.Range('E20').Value = '  -0.25%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '511.64'
$ws.Range('E21').Value = '  +3.40%  '

# Row 22
$ws.Range('E22').Value = '  -3.34%  '

# Row 23
$ws.Range('E23').Value = '  -3.86%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.19'
$ws.Range('E24').Value = '  -2.50%  '

# Row 25
$ws.Range('E25').Value = '  -4.92%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.47'
$ws.Range('E26').Value = '  +3.12%  '

# Row 27
$ws.Range('D27').Value = '3.914.68'
$ws.Range('E27').Value = '  -1.69%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.11'
$ws.Range('E28').Value = '  -4.91%  '

# Row 30
$ws.Range('E30').Value = '  -7.52%  '

# Row 31
$ws.Range('E31').Value = '  -6.23%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.24'
$ws.Range('E33').Value = '  -2.63%  '

# Row 34
$ws.Range('E34').Value = '  -1.94%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.18'
$ws.Range('E35').Value = '  -0.14%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.167'
$ws.Range('E36').Value = '  -5.27%  '

# Row 38
$ws.Range('E38').Value = '  -1.58%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.49'
$ws.Range('E39').Value = '  +8.47%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0989'
$ws.Range('E40').Value = '  -3.51%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.29'
$ws.Range('E41').Value = '  -1.62%  '

# Row 42
$ws.Range('E42').Value = '  -3.82%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  -0.17%  '

# Row 45
$ws.Range('E45').Value = '  -6.89%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '160.74'
$ws.Range('E46').Value = '  -1.64%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '48.86'
$ws.Range('E47').Value = '  -0.05%  '

# Row 48
$ws.Range('E48').Value = '  -5.37%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.54'
$ws.Range('E49').Value = '  -2.74%  '

# Row 50
$ws.Range('E50').Value = '  -0.76%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '406.50'
$ws.Range('E51').Value = '  -5.20%  '
